$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 303
    3  = 396
    4  = 390
    5  = 296
    6  = 190
    7  = 155
    8  = 158
    9  = 176
    10 = 166
    11 = 154
    12 = 153
    13 = 166
    14 = 174
    15 = 167
    16 = 153
    17 = 156
    18 = 165
    19 = 167
    20 = 147
    21 = 124
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
